$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    39.661576838990712,
    295.76621513053516,
    74.733977165428001,
    126.88069883751206,
    130.20867875123534,
    430.88892641510722,
    254.26403544915837,
    94.801936647429827,
    143.23278033753741,
    117.75149697655652,
    258.60655480161415,
    52.443698438959167,
    27.074246049955519,
    257.87275622842628,
    33.04314485073985,
    173.33989559288116,
    208.17209941609539,
    757.10632321431433,
    163.14034917521431,
    89.559845473220037,
    128.72229685243937,
    128.43969387593808,
    225.89213287039405,
    82.947213971479798,
    113.44491440584073,
    5358.6679975083171,
    58.803254715400463,
    260.54709154766317,
    203.6717821177281,
    265.07027376910264,
    345.49898453816445,
    104.6751717930136,
    68.329405730533367,
    160.5091910730176,
    163.86053316942341,
    531.99168167310359,
    46.942016262882142,
    159.23440579921328,
    113.72894733620301,
    744.04509883765877,
    224.41242259378564,
    268.72760405871361,
    135.25711848068877,
    187.89376781094867,
    233.5135551461384,
    51.13615847776105,
    242.75865789689698,
    361.93790454455956,
    204.24147580917403,
    136.86445181394473,
    398.70533643903451,
    35.837263721126227,
    476.16138808071042,
    389.02315538207284,
    352.51787408860872,
    572.3550149539177,
    83.052129353280222,
    151.18084098533325,
    323.08532054171104,
    185.2535526837807,
    106.09015401974783,
    185.00269291138164,
    85.741034010410914,
    652.21707044179755,
    373.67707282607932,
    14809.972877809972,
    133.82224093694566,
    251.80782104352272,
    133.06565696596417,
    457.52493923349601,
    134.84257473493412,
    115.27995318571632,
    3571.4904413068298,
    1467.4656695062847,
    1136.849486141376,
    526.75866053598588,
    75.165634690516313,
    63.283667760863146,
    2729.143040590036,
    14824.323764916691,
    695.62553332888876,
    99.528070339120518,
    107.16072071101756,
    66.204617619713545,
    169.54592441811812,
    1855.3326954348984,
    86.891437461189696,
    66.915592483345065,
    76.748618130886641,
    2889.9817903173771,
    6148.1391984159764,
    2294.2269463621556,
    51.869931312899503,
    180.69882251811921,
    120.24099369651071,
    1665.2022897403235,
    92.96112726645039,
    1434.4389208812911,
    2001.7734309749542,
    47.640038874502579,
    1138.5603361865105,
    127.60089229038861,
    1690.208403760899
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 3).Value = $values[$i]
}

$ws.Range("A1:C103").Select()
